# 🔄 MAJ automatique BRVM via GitHub Actions
#
# Refreshes the "Recommandations" sheet (now 36 titles instead of 30,
# rows 2..37) and the "Top_YTD" sheet (still 10 sector rows, 2..11) with
# the latest computed figures, re-sorted by performance.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Recommandations
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Recommandations")

$reco = @(
    @("BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 4, 652.53, 159.95, "🟡 Observer", "➖ Neutre"),
    @("BRVM - SERVICES FINANCIERS", 0, 4, 582.91, 146.61, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRESTIGE", 0, 4, 571.78, 143.73, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIELS", 0, 4, 527.4, 133.71, "🟡 Observer", "➖ Neutre"),
    @("BRVM - ENERGIE", 0, 4, 442.87, 111.83, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DE BASE    (**)", 0, 2, 431.6, 215.68, "🟡 Observer", "➖ Neutre"),
    @("BRVM - SERVICES PUBLICS", 0, 4, 427.07, 107.35, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TELECOMMUNICATIONS", 0, 4, 373.13, 93.52, "🟡 Observer", "➖ Neutre"),
    @("BRVM – COMPOSITE TOTAL RETURN    (**)", 0, 2, 263.52, 130.87, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL    (**)", 0, 1, 217.75, 217.75, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL     (**)", 0, 1, 216.13, 216.13, "🟡 Observer", "➖ Neutre"),
    @("EVIOSYS PACKAGING SIEM CI (SEMC)", 4, 0, 28.86, 6.94, "🟢 Achat", "✅ Renforcer"),
    @("VIVO ENERGY CI (SHEC)", 2, 0, 9.69, 4.73, "🟡 Observer", "➖ Neutre"),
    @("ORAGROUP TOGO (ORGT)", 2, 0, 9.58, 6.04, "🟡 Observer", "➖ Neutre"),
    @("SICOR CI (SICC)", 1, 0, 7.4, 7.4, "🟡 Observer", "➖ Neutre"),
    @("UNILEVER CI (UNLC)", 1, 0, 7.14, 7.14, "🟡 Observer", "➖ Neutre"),
    @("ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)", 2, 1, 6.64, -1.77, "🟡 Observer", "👀 À surveiller"),
    @("BANK OF AFRICA BF (BOABF)", 1, 0, 5.07, 5.07, "🟡 Observer", "➖ Neutre"),
    @("SOCIETE IVOIRIENNE DE BANQUE  (SIBC)", 1, 0, 4.1, 4.1, "🟡 Observer", "➖ Neutre"),
    @("NESTLE CI (NTLC)", 1, 0, 3.24, 3.24, "🟡 Observer", "➖ Neutre"),
    @("SOGB CI (SOGC)", 1, 0, 3.24, 3.24, "🟡 Observer", "➖ Neutre"),
    @("ONATEL BF (ONTBF)", 1, 0, 1.22, 1.22, "🟡 Observer", "➖ Neutre"),
    @("SERVAIR ABIDJAN CI (ABJC)", 0, 1, -1.8, -1.8, "🟡 Observer", "➖ Neutre"),
    @("UNIWAX CI (UNXC)", 0, 1, -2.13, -2.13, "🟡 Observer", "➖ Neutre"),
    @("LOTERIE NATIONALE DU BENIN (LNBB)", 0, 1, -2.21, -2.21, "🟡 Observer", "➖ Neutre"),
    @("SICABLE CI (CABC)", 0, 1, -2.34, -2.34, "🟡 Observer", "➖ Neutre"),
    @("SMB CI (SMBC)", 0, 1, -3.16, -3.16, "🟡 Observer", "➖ Neutre"),
    @("BERNABE CI (BNBC)", 0, 1, -3.19, -3.19, "🟡 Observer", "➖ Neutre"),
    @("SETAO CI (STAC)", 0, 1, -3.4, -3.4, "🟡 Observer", "➖ Neutre"),
    @("SOLIBRA CI (SLBC)", 1, 1, -3.45, -4.5, "🟡 Observer", "👀 À surveiller"),
    @("TOTALENERGIES MARKETING SN (TTLS)", 1, 1, -3.92, 3.19, "🟡 Observer", "👀 À surveiller"),
    @("SAFCA CI (SAFC)", 0, 1, -3.93, -3.93, "🟡 Observer", "➖ Neutre"),
    @("ECOBANK TRANS. INCORP. TG (ETIT)", 0, 1, -4.35, -4.35, "🟡 Observer", "➖ Neutre"),
    @("FILTISAC CI (FTSC)", 1, 2, -5.42, 7.27, "🟡 Observer", "👀 À surveiller"),
    @("CFAO MOTORS CI (CFAC)", 0, 3, -13.11, -4.21, "🔴 Vente", "⚠️ Risque de décrochage"),
    @("NEI-CEDA CI (NEIC)", 0, 3, -17.62, -2.91, "🔴 Vente", "⚠️ Risque de décrochage")
)

$r = 2
foreach ($row in $reco) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
    $r++
}

# ---------------------------------------------------------------------
# Sheet 2: Top_YTD
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Top_YTD")

$ytd = @(
    @("BRVM - CONSOMMATION DISCRETIONNAIRE", 4691.3),
    @("BRVM - SERVICES FINANCIERS", 3545.93),
    @("BRVM-PRESTIGE", 3383.53),
    @("BRVM - INDUSTRIELS", 2789.36),
    @("BRVM - ENERGIE", 1871.42),
    @("BRVM - SERVICES PUBLICS", 1727.78),
    @("BRVM - TELECOMMUNICATIONS", 1295.59),
    @("BRVM - CONSOMMATION DE BASE    (**)", 897.3),
    @("BRVM – COMPOSITE TOTAL RETURN    (**)", 437.12),
    @("BRVM-PRINCIPAL    (**)", 217.75)
)

$r = 2
foreach ($row in $ytd) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $r++
}
